# Append the latest daily profit figure (run on 2025-10-21) as a new row
# at the bottom of the existing Date/Profit table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the current data block.
$lastRow = $ws.UsedRange.Rows.Count
$newRow  = $lastRow + 1

$dateCell   = $ws.Cells.Item($newRow, 1)
$profitCell = $ws.Cells.Item($newRow, 2)

# Force the date column to be stored as plain text (matching the existing
# rows, which hold the date as literal "MM/DD/YYYY" text rather than a
# serial date value), then restore the default "Normal" style so no stray
# number-format/style is left behind on the cell.
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/21/2025"
$dateCell.Style = "Normal"

$profitCell.Value = 9786.5
